$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new columns before column D (old D:K shifts to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: copy number formatting from column F (the shifted former column D)
# onto the two new columns D:E so the new cells match the format of their row
# (s="2" date format on header rows, s="3" comma-number format elsewhere).
# Only the three data blocks (income statement / balance sheet / cash flow)
# actually have cells in D:K - the section-title rows in between (37, 79) must
# stay untouched, so the format copy is done per block instead of D7:E102.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 3: populate the new D:E columns with the newest-quarter figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43372
$ws.Range("D8").Value = 4942400
$ws.Range("E8").Value = 4891600
$ws.Range("D9").Value = 4397000
$ws.Range("E9").Value = 4347500
$ws.Range("D10").Value = 545400
$ws.Range("E10").Value = 544100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 39200
$ws.Range("E14").Value = 19800
$ws.Range("D15").Value = 12500
$ws.Range("E15").Value = 12700
$ws.Range("D17").Value = 4603000
$ws.Range("E17").Value = 4541500
$ws.Range("D18").Value = 339400
$ws.Range("E18").Value = 350100
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 462000
$ws.Range("E21").Value = 469900
$ws.Range("D22").Value = 21300
$ws.Range("E22").Value = 21200
$ws.Range("D23").Value = 318100
$ws.Range("E23").Value = 328900
$ws.Range("D24").Value = 74900
$ws.Range("E24").Value = 64500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 243200
$ws.Range("E26").Value = 264400
$ws.Range("D27").Value = 216200
$ws.Range("E27").Value = 245600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -4000
$ws.Range("E29").Value = 6900
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 212200
$ws.Range("E33").Value = 252500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 212200
$ws.Range("E35").Value = 252500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43372
$ws.Range("D41").Value = 1493200
$ws.Range("E41").Value = 1198600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 2880300
$ws.Range("E43").Value = 3301800
$ws.Range("D44").Value = 1196800
$ws.Range("E44").Value = 1285200
$ws.Range("D45").Value = 710200
$ws.Range("E45").Value = 782100
$ws.Range("D46").Value = 6280500
$ws.Range("E46").Value = 6567700
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2598100
$ws.Range("E48").Value = 2536200
$ws.Range("D49").Value = 1801800
$ws.Range("E49").Value = 1409400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 920300
$ws.Range("E52").Value = 1496100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 11600700
$ws.Range("E54").Value = 12009400
$ws.Range("D57").Value = 2862800
$ws.Range("E57").Value = 3041800
$ws.Range("D58").Value = 22800
$ws.Range("E58").Value = 17200
$ws.Range("D59").Value = 1615000
$ws.Range("E59").Value = 1720900
$ws.Range("D60").Value = 4500600
$ws.Range("E60").Value = 4779900
$ws.Range("D61").Value = 1941000
$ws.Range("E61").Value = 1946600
$ws.Range("D62").Value = 640400
$ws.Range("E62").Value = 685600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 7400000
$ws.Range("E66").Value = 7716500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 4113600
$ws.Range("E72").Value = 4953800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 4200700
$ws.Range("E76").Value = 4292900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43372
$ws.Range("D81").Value = 212200
$ws.Range("E81").Value = 252500
$ws.Range("D83").Value = 122600
$ws.Range("E83").Value = 119800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 758200
$ws.Range("E89").Value = 267900
$ws.Range("D91").Value = -184300
$ws.Range("E91").Value = -160500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -182300
$ws.Range("E94").Value = -153900
$ws.Range("D96").Value = -44200
$ws.Range("E96").Value = -45400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -277300
$ws.Range("E100").Value = -246100
$ws.Range("D101").Value = -3700
$ws.Range("E101").Value = -9800
$ws.Range("D102").Value = 294900
$ws.Range("E102").Value = -141900
